{"js": "// 1. \"As a Presales Architect\" -> \"As a Lead Presales Architect\"\nconst body = context.document.body;\nconst searchResults = body.search(\"As a Presales Architect\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"As a Lead Presales Architect\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Insert a new bullet paragraph \"Integration of technical architects within the\n//    dynamic Sales enviornment\" right after \"Optimum sizing as per customer use-case\n//    and budget\" and before \"Providing post-sales support to newly onboarded Customers.\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Optimum sizing as per customer use-case and budget\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (anchor) {\n  anchor.insertParagraph(\n    \"Integration of technical architects within the dynamic Sales enviornment\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"As a Presales Architect\" -> \"As a Lead Presales Architect\"\n$find = $d.Content.Find\n$find.Text = \"As a Presales Architect\"\n$find.Replacement.Text = \"As a Lead Presales Architect\"\n$find.Execute(\n    [ref]\"As a Presales Architect\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"As a Lead Presales Architect\", 2\n)\n\n# 2. Insert a new bullet paragraph \"Integration of technical architects within the\n#    dynamic Sales enviornment\" right after \"Optimum sizing as per customer use-case\n#    and budget\" and before \"Providing post-sales support to newly onboarded Customers.\"\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"Optimum sizing as per customer use-case and budget\") {\n        $p.Range.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newPara.Range.Text = \"Integration of technical architects within the dynamic Sales enviornment\"\n        break\n    }\n}\n"}
